# Apply "repull data, push all data, mean calculation" updates to column F (dSF)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F
$updates = @{
    3  = -4
    6  = 4
    10 = -3
    12 = 6
    13 = 1
    17 = 5
    19 = 6
    24 = 1
    27 = -6
    36 = -7
    37 = -3
    38 = -6
    41 = -6
    45 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
